# Update dSF (column F) values on Sheet1 per the repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = -6
$ws.Range("F4").Value = -1
$ws.Range("F8").Value = -2
$ws.Range("F10").Value = 4
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = -2
$ws.Range("F23").Value = -4
$ws.Range("F30").Value = -3
$ws.Range("F32").Value = 0
